$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 705.1667
$ws.Range("I19").Value = 683.5833
$ws.Range("J19").Value = 748.3333
$ws.Range("K19").Value = 683.5833
$ws.Range("L19").Value = 748.3333
$ws.Range("M19").Value = -508.5833
$ws.Range("N19").Value = -1098.3333

$ws.Range("H86").Value = 1806.4706
$ws.Range("I86").Value = 1507.6154
$ws.Range("J86").Value = 2777.75
$ws.Range("K86").Value = 1507.6154
$ws.Range("L86").Value = 2777.75
$ws.Range("M86").Value = -384.6153999999999
$ws.Range("N86").Value = -5023.75

$ws.Range("H89").Value = 1806.4706
$ws.Range("I89").Value = 1507.6154
$ws.Range("J89").Value = 2777.75
$ws.Range("K89").Value = 7538.076999999999
$ws.Range("L89").Value = 13888.75
$ws.Range("M89").Value = -1922.076999999999
$ws.Range("N89").Value = -25120.75

$ws.Range("H92").Value = 2050.1875
$ws.Range("I92").Value = 1475
$ws.Range("J92").Value = 3315.6
$ws.Range("K92").Value = 1475
$ws.Range("L92").Value = 3315.6
$ws.Range("M92").Value = -227
$ws.Range("N92").Value = -5811.6

$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 2000
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -1549
$ws.Range("N94").Value = -2902

$ws.Range("H111").Value = 6900
$ws.Range("I111").Value = 6900
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 20700
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = $null
$ws.Range("N111").Value = -17633

$ws.Range("H113").Value = 5800
$ws.Range("I113").Value = 5800
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5800
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = -2546

$ws.Range("H125").Value = 2052.6428
$ws.Range("J125").Value = 3577.7144
$ws.Range("L125").Value = 32199.4296
$ws.Range("N125").Value = -37119.4296

$ws.Range("H129").Value = 814.6724
$ws.Range("I129").Value = 376.14285
$ws.Range("J129").Value = 954.2045000000001
$ws.Range("K129").Value = 1128.42855
$ws.Range("L129").Value = 2862.6135
$ws.Range("M129").Value = 3871.57145
$ws.Range("N129").Value = -12862.6135

$ws.Range("H138").Value = 2196905.5
$ws.Range("I138").Value = 2034.8518
$ws.Range("J138").Value = 3406324
$ws.Range("K138").Value = 6104.555399999999
$ws.Range("L138").Value = 10218972
$ws.Range("M138").Value = -964.5553999999993
$ws.Range("N138").Value = -10229252

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224

$ws.Range("H32").Value = 4237215
$ws.Range("I32").Value = 4711241.5
$ws.Range("J32").Value = 30227.625
$ws.Range("K32").Value = 4711241.5
$ws.Range("L32").Value = 30227.625
$ws.Range("M32").Value = -4710954.5
$ws.Range("N32").Value = -30801.625

$ws.Range("H74").Value = 4755561
$ws.Range("I74").Value = 6784441
$ws.Range("J74").Value = 63776
$ws.Range("K74").Value = 6784441
$ws.Range("L74").Value = 63776
$ws.Range("M74").Value = -6783567
$ws.Range("N74").Value = -65524

$ws.Range("H77").Value = 4755561
$ws.Range("I77").Value = 6784441
$ws.Range("J77").Value = 63776
$ws.Range("K77").Value = 33922205
$ws.Range("L77").Value = 318880
$ws.Range("M77").Value = -33917837
$ws.Range("N77").Value = -327616

$ws.Range("H110").Value = 556470.6
$ws.Range("I110").Value = 556470.6
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 556470.6
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = $null
$ws.Range("N110").Value = -554425.6

$ws.Range("H132").Value = 73476.64
$ws.Range("I132").Value = 49164.24
$ws.Range("J132").Value = 146413.86
$ws.Range("K132").Value = 147492.72
$ws.Range("L132").Value = 439241.58
$ws.Range("M132").Value = -144962.72
$ws.Range("N132").Value = -444301.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1230

$ws.Range("H134").Value = 3333.4055
$ws.Range("I134").Value = 2677.4546
$ws.Range("K134").Value = 8032.3638
$ws.Range("M134").Value = -5497.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8312.588
$ws.Range("I31").Value = 15355.737
$ws.Range("J31").Value = 1940.2142
$ws.Range("K31").Value = 15355.737
$ws.Range("L31").Value = 1940.2142
$ws.Range("M31").Value = -15060.737
$ws.Range("N31").Value = -2530.2142

$ws.Range("H34").Value = 8312.588
$ws.Range("I34").Value = 15355.737
$ws.Range("J34").Value = 1940.2142
$ws.Range("K34").Value = 15355.737
$ws.Range("L34").Value = 1940.2142
$ws.Range("M34").Value = -15153.737
$ws.Range("N34").Value = -2344.2142

$ws.Range("H105").Value = 1260
$ws.Range("I105").Value = 1286.7894
$ws.Range("K105").Value = 1286.7894
$ws.Range("M105").Value = 460.2106000000001

$ws.Range("H132").Value = 26610.293
$ws.Range("I132").Value = 1734.7084
$ws.Range("J132").Value = 61728.766
$ws.Range("K132").Value = 5204.1252
$ws.Range("L132").Value = 185186.298
$ws.Range("M132").Value = -2674.1252
$ws.Range("N132").Value = -190246.298

$ws.Range("H133").Value = 40031.7
$ws.Range("J133").Value = 40031.7
$ws.Range("L133").Value = 40031.7
$ws.Range("N133").Value = -45091.7

$ws.Range("H134").Value = 23472.389
$ws.Range("I134").Value = 1401.0646
$ws.Range("J134").Value = 61484.11
$ws.Range("K134").Value = 4203.1938
$ws.Range("L134").Value = 184452.33
$ws.Range("M134").Value = -1668.1938
$ws.Range("N134").Value = -189522.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 986.6818
$ws.Range("I14").Value = 986.6818
$ws.Range("K14").Value = 2960.0454
$ws.Range("M14").Value = -2787.0454

$ws.Range("H131").Value = 857.61365
$ws.Range("J131").Value = 913.6486
$ws.Range("L131").Value = 2740.9458
$ws.Range("N131").Value = -12820.9458

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2880.2
$ws.Range("I122").Value = 2600.25
$ws.Range("K122").Value = 7800.75
$ws.Range("M122").Value = -5350.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 12627.728
$ws.Range("I16").Value = 9000.143
$ws.Range("J16").Value = 18976
$ws.Range("K16").Value = 9000.143
$ws.Range("L16").Value = 18976
$ws.Range("M16").Value = -8830.143
$ws.Range("N16").Value = -19316

$ws.Range("H42").Value = 7800
$ws.Range("J42").Value = 7800
$ws.Range("L42").Value = 7800
$ws.Range("N42").Value = -8926

$ws.Range("H46").Value = 920.6923
$ws.Range("I46").Value = 498
$ws.Range("K46").Value = 498
$ws.Range("M46").Value = -310

$ws.Range("H49").Value = 7800
$ws.Range("J49").Value = 7800
$ws.Range("L49").Value = 7800
$ws.Range("N49").Value = -8094

$ws.Range("H122").Value = 2444.4
$ws.Range("I122").Value = 2294.6667
$ws.Range("J122").Value = 2669
$ws.Range("K122").Value = 6884.000100000001
$ws.Range("L122").Value = 8007
$ws.Range("M122").Value = -4434.000100000001
$ws.Range("N122").Value = -12907

$ws.Range("H132").Value = 38610.5
$ws.Range("I132").Value = 2197.0527
$ws.Range("J132").Value = 115483.336
$ws.Range("K132").Value = 6591.158100000001
$ws.Range("L132").Value = 346450.008
$ws.Range("M132").Value = -4061.158100000001
$ws.Range("N132").Value = -351510.008

$ws.Range("H136").Value = 59799.46
$ws.Range("I136").Value = 27099.309
$ws.Range("J136").Value = 175736.36
$ws.Range("K136").Value = 81297.927
$ws.Range("L136").Value = 527209.08
$ws.Range("M136").Value = -78747.927
$ws.Range("N136").Value = -532309.08

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 10164.143
$ws.Range("I49").Value = 5000
$ws.Range("J49").Value = 11024.833
$ws.Range("K49").Value = 5000
$ws.Range("L49").Value = 11024.833
$ws.Range("M49").Value = -4770
$ws.Range("N49").Value = -11484.833

$ws.Range("H107").Value = 666.6667
$ws.Range("I107").Value = 500
$ws.Range("K107").Value = 1500
$ws.Range("M107").Value = 420

$ws.Range("H132").Value = 59179.43
$ws.Range("I132").Value = 43467
$ws.Range("J132").Value = 93461.09
$ws.Range("K132").Value = 130401
$ws.Range("L132").Value = 280383.27
$ws.Range("M132").Value = -127871
$ws.Range("N132").Value = -285443.27
